# The "raw" sheet used to carry a single stray value at AC2 (=15) plus the
# padding placeholder at F11. The data-generation script was switched from
# openpyxl to numpy, which now fills a full A2:A11 column of 15s instead of
# the lone AC2 cell, and leaves the view scrolled/selected further right
# (over near the W/S columns) instead of parked at F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw")
$ws.Activate()

# Drop the old stray value that openpyxl had written at AC2.
$ws.Range("AC2").ClearContents()

# numpy-driven fill: A2:A11 each get the value 15.
$ws.Range("A2:A11").Value = 15

# Scroll the window so column S is the left-most visible column, then settle
# the selection on W10 (previously the view was parked at F11).
$excel.ActiveWindow.ScrollColumn = 19
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("W10").Select()
